$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.672.94'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '1.594.94'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.01'
$ws.Range('E5').Value = '  -0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.513'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  -1.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.47'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').Value = '1.820.19'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = '1.599.35'
$ws.Range('E14').Value = '  -0.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.69'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = '26.660.75'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.00'
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '207.66'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.79'
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.33'
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.87'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.58'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.20'
$ws.Range('E27').Value = '  -2.33%  '
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.23'
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0505'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.15'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  -0.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.662'
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('E34').Value = '  +0.29%  '
$ws.Range('D35').Value = '1.285.79'
$ws.Range('E35').Value = '  -3.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.45'
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0171'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.841'
$ws.Range('E39').Value = '  +1.98%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('E41').Value = '  +1.39%  '
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.785'
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.58'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('D45').Value = '1.732.11'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.903'
$ws.Range('E46').Value = '  +8.88%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.08'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  +0.19%  '
